$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A12").Value = '{''Gender'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A13").Value = '{''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A14").Value = '{''Gender'': ''1'', ''HDI'': ''1''}'
$ws.Range("A15").Value = '{''Gender'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A17").Value = '{''Gender'': ''1'', ''Student'': ''1''}'
$ws.Range("A18").Value = '{''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A19").Value = '{''Gender'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A20").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A21").Value = '{''Gender'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A22").Value = '{''Gender'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A23").Value = '{''Gender'': ''1'', ''Hobby'': ''1''}'
$ws.Range("A24").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A26").Value = '{''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A27").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A28").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A29").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A30").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A31").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A32").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A33").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A34").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A35").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A36").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A37").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''Gender'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A39").Value = '{''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A40").Value = '{''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A41").Value = '{''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A42").Value = '{''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A43").Value = '{''Gender'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A44").Value = '{''Gender'': ''1'', ''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A45").Value = '{''Dependents'': ''2'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A46").Value = '{''Gender'': ''1'', ''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A47").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A48").Value = '{''Gender'': ''1'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A49").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A50").Value = '{''Gender'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A51").Value = '{''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A53").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1''}'
$ws.Range("A55").Value = '{''Gender'': ''1'', ''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A56").Value = '{''Gender'': ''1'', ''DevType'': ''2''}'
$ws.Range("A58").Value = '{''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A59").Value = '{''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A60").Value = '{''Gender'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A61").Value = '{''Gender'': ''1'', ''HoursComputer'': ''2''}'
$ws.Range("A62").Value = '{''HoursComputer'': ''2'', ''SexualOrientation'': ''1''}'
